# Update absenteeism data rows 2-11 with new values per the source diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @{ Row = 2;  A = 89668; B = "Alice da Cunha";          C = "Engenharia";             D = "Viagem de negócios"; E = 2; F = 45078; G = 5118.18 }
    @{ Row = 3;  A = 38876; B = "Samuel Porto";             C = "Recursos Humanos";       D = "Outros";              E = 7; F = 45094; G = 6696.52 }
    @{ Row = 4;  A = 91380; B = "Sra. Ana Clara Costela";   C = "Atendimento ao Cliente"; D = "Problemas pessoais";  E = 3; F = 45088; G = 5429.68 }
    @{ Row = 5;  A = 71062; B = "Diogo Nunes";               C = "Vendas";                 D = "Outros";              E = 5; F = 45100; G = 10813.33 }
    @{ Row = 6;  A = 57991; B = "Augusto da Rocha";          C = "Recursos Humanos";       D = "Problemas pessoais";  E = 3; F = 45094; G = 8661.82 }
    @{ Row = 7;  A = 91064; B = "Luana da Mata";             C = "TI";                     D = "Viagem de negócios"; E = 4; F = 45086; G = 7389.34 }
    @{ Row = 8;  A = 4197;  B = "Sra. Stella Monteiro";      C = "Atendimento ao Cliente"; D = "Doença";              E = 1; F = 45083; G = 7813.02 }
    @{ Row = 9;  A = 38131; B = "Luiz Henrique Correia";     C = "Vendas";                 D = "Outros";              E = 7; F = 45085; G = 3632.76 }
    @{ Row = 10; A = 50052; B = "Sra. Ana Julia Cardoso";    C = "Engenharia";             D = "Viagem de negócios"; E = 7; F = 45087; G = 5334.9 }
    @{ Row = 11; A = 36114; B = "Diogo Almeida";             C = "Atendimento ao Cliente"; D = "Doença";              E = 3; F = 45091; G = 11017.43 }
)

foreach ($rowData in $data) {
    $r = $rowData.Row
    $ws.Cells.Item($r, 1).Value = $rowData.A
    $ws.Cells.Item($r, 2).Value = $rowData.B
    $ws.Cells.Item($r, 3).Value = $rowData.C
    $ws.Cells.Item($r, 4).Value = $rowData.D
    $ws.Cells.Item($r, 5).Value = $rowData.E
    $ws.Cells.Item($r, 6).Value = $rowData.F
    $ws.Cells.Item($r, 7).Value = $rowData.G
}
